$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddCoverageTeam")

# Update the coverage-team row values (refreshed after Db refresh)
$ws.Range("C2").Value = "Primary"
$ws.Range("D2").Value = "Education"

# Reflect the new active selection on this sheet
$ws.Activate()
$ws.Range("E8").Select()
